$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "fantasy points" column (G),
# shifting it to I, to make room for "height" and "weight".
$ws.Range("G:H").EntireColumn.Insert()

# Header row
$ws.Range("G1").Value = "height"
$ws.Range("H1").Value = "weight"

# Copy the header style (bold + border) from column F's header onto the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats

$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 6    # G: height
    $ws.Cells.Item($r, 8).Value = 209  # H: weight
}
